# feat: add 2022-Q4 data
#
# - Insert a new "2022-Q4" worksheet right after "总计" (becomes the 2nd tab),
#   populated with the new quarter's fund-holding data.
# - Update the "总计" (totals) sheet: insert a new top data row for 2022-Q4
#   and renumber the existing index column.
# - "2022-Q3", "2022-Q2" and "2022-Q1" sheets keep their data untouched; they
#   simply shift right in tab order because the new sheet is inserted before them.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Create the new "2022-Q4" sheet, positioned right after "总计".
# ---------------------------------------------------------------------------
$zj = $wb.Worksheets.Item("总计")

$q4 = $wb.Worksheets.Add([System.Type]::Missing, $zj)
$q4.Name = "2022-Q4"

# NOTE: worksheet references must be (re)fetched *after* Worksheets.Add(),
# otherwise handles obtained beforehand paste/copy incorrectly.
$q3 = $wb.Worksheets.Item("2022-Q3")

# Borrow the header/body cell formatting from the existing "2022-Q3" sheet so
# the new sheet matches the look (bold centered bordered header row + index
# column) of all the other quarterly sheets.
$q3.Range("A1:H3").Copy()
$q4.Range("A1").PasteSpecial(-4122)  # xlPasteFormats

$q4.Range("A3:H3").Copy()
$q4.Range("A4:H7").PasteSpecial(-4122)  # xlPasteFormats

# Header row
$q4.Range("B1").Value = "基金代码"
$q4.Range("C1").Value = "基金名称"
$q4.Range("D1").Value = "基金规模"
$q4.Range("E1").Value = "股票总仓位"
$q4.Range("F1").Value = "仓位占比"
$q4.Range("G1").Value = "持有市值(亿元)"
$q4.Range("H1").Value = "仓位排名"

# The fund-code / size / position columns are stored as text even though they
# look numeric (e.g. leading-zero fund codes, "1.86"), so force text format
# before writing those values.
$q4.Range("B2:B7").NumberFormat = "@"
$q4.Range("D2:G7").NumberFormat = "@"

$q4.Range("A2").Value = 0
$q4.Range("B2").Value = "004138"
$q4.Range("C2").Value = "上银鑫达灵活配置混合A"
$q4.Range("D2").Value = "1.86"
$q4.Range("E2").Value = "80.66"
$q4.Range("F2").Value = "3.01"
$q4.Range("G2").Value = "0.0560"
$q4.Range("H2").Value = 9

$q4.Range("A3").Value = 1
$q4.Range("B3").Value = "016285"
$q4.Range("C3").Value = "汇丰晋信龙头优势混合A"
$q4.Range("D3").Value = "1.61"
$q4.Range("E3").Value = "58.20"
$q4.Range("F3").Value = "2.60"
$q4.Range("G3").Value = "0.0419"
$q4.Range("H3").Value = 5

$q4.Range("A4").Value = 2
$q4.Range("B4").Value = "012334"
$q4.Range("C4").Value = "上银慧尚6个月持有期混合A"
$q4.Range("D4").Value = "2.99"
$q4.Range("E4").Value = "22.23"
$q4.Range("F4").Value = "0.76"
$q4.Range("G4").Value = "0.0227"
$q4.Range("H4").Value = 8

$q4.Range("A5").Value = 3
$q4.Range("B5").Value = "016286"
$q4.Range("C5").Value = "汇丰晋信龙头优势混合C"
$q4.Range("D5").Value = "0.26"
$q4.Range("E5").Value = "58.20"
$q4.Range("F5").Value = "2.60"
$q4.Range("G5").Value = "0.0068"
$q4.Range("H5").Value = 5

$q4.Range("A6").Value = 4
$q4.Range("B6").Value = "012335"
$q4.Range("C6").Value = "上银慧尚6个月持有期混合C"
$q4.Range("D6").Value = "0.14"
$q4.Range("E6").Value = "22.23"
$q4.Range("F6").Value = "0.76"
$q4.Range("G6").Value = "0.0011"
$q4.Range("H6").Value = 8

$q4.Range("A7").Value = 5
$q4.Range("B7").Value = "015753"
$q4.Range("C7").Value = "上银鑫达灵活配置混合C"
$q4.Range("D7").Value = "0.00"
$q4.Range("E7").Value = "80.66"
$q4.Range("F7").Value = "3.01"
# G7 is a genuine numeric 0 (not text), unlike the other G-column cells.
$q4.Range("G7").NumberFormat = "General"
$q4.Range("G7").Value = 0
$q4.Range("H7").Value = 9

# Match the page margins used by the sibling quarterly sheets.
$q4.PageSetup.LeftMargin = 54
$q4.PageSetup.RightMargin = 54
$q4.PageSetup.TopMargin = 72
$q4.PageSetup.BottomMargin = 72
$q4.PageSetup.HeaderMargin = 36
$q4.PageSetup.FooterMargin = 36

# ---------------------------------------------------------------------------
# 2. Update "总计": add the 2022-Q4 total row at the top, renumber the rest.
# ---------------------------------------------------------------------------
$zj.Range("A2").Value = 0
$zj.Range("B2").Value = "2022-Q4"
$zj.Range("C2").Value = 6
$zj.Range("D2").Value = 0.13

$zj.Range("A3").Value = 1
$zj.Range("B3").Value = "2022-Q3"
$zj.Range("C3").Value = 2
$zj.Range("D3").Value = 0.05

$zj.Range("A4").Value = 2
$zj.Range("B4").Value = "2022-Q2"
$zj.Range("C4").Value = 2
$zj.Range("D4").Value = 0.21

$zj.Range("A5").Value = 3
$zj.Range("B5").Value = "2022-Q1"
$zj.Range("C5").Value = 4
$zj.Range("D5").Value = 0.35

# A5 is a brand-new cell; copy the index-column style down from A4 so it
# keeps the same bold/centered/bordered look as the rest of column A.
$zj.Range("A4").Copy()
$zj.Range("A5").PasteSpecial(-4122)  # xlPasteFormats
$zj.Range("A5").Value = 3

# ---------------------------------------------------------------------------
# 3. Restore "2022-Q1" as the active/selected sheet (it was the active tab
#    before the edit and nothing in this change should move the user's
#    selection away from it).
# ---------------------------------------------------------------------------
$q1 = $wb.Worksheets.Item("2022-Q1")
$q1.Activate()
$q1.Range("A1").Select()
